$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AID (col B) and Intake Date (col E) are stored as plain TEXT in this sheet
# (not real numbers/dates). Force text entry so Excel's smart type detection
# doesn't turn "58710884" into a number or "06/14/2025" into a date serial,
# then clear the temporary number format back off so no extra formatting is
# left behind on the cells.
$ws.Range("B2:B9").NumberFormat = "@"
$ws.Range("E2:E9").NumberFormat = "@"

# --- Row 2: Livia -> SMORES (Cat, Foster Home) ---
$ws.Range("A2").Value = "SMORES"
$ws.Range("B2").Value = "58710884"
$ws.Range("C2").Value = "Cat"
$ws.Range("D2").Value = "Foster Home"
$ws.Range("E2").Value = "06/14/2025"
$ws.Range("F2").Value = "Domestic Shorthair"
$ws.Range("G2").Value = "2y 2m 5d"
$ws.Range("H2").Value = "In Foster"
$ws.Range("I2").Value = 65.8

# --- Row 3: Brody -> Petey (Dog Adoptions D) ---
$ws.Range("A3").Value = "Petey"
$ws.Range("B3").Value = "58804353"
$ws.Range("C3").Value = "Dog"
$ws.Range("D3").Value = "Dog Adoptions D"
$ws.Range("E3").Value = "07/29/2025"
$ws.Range("F3").Value = "Mixed Breed, Large (over 44 lbs fully grown)"
$ws.Range("G3").Value = "11m 22d"
$ws.Range("H3").Value = "Available"
$ws.Range("I3").Value = 20.9

# --- Row 4: Karma -> Remy (Dog Adoptions D) ---
$ws.Range("A4").Value = "Remy"
$ws.Range("B4").Value = "58959672"
$ws.Range("C4").Value = "Dog"
$ws.Range("D4").Value = "Dog Adoptions D"
$ws.Range("E4").Value = "07/31/2025"
$ws.Range("F4").Value = "Mixed Breed, Large (over 44 lbs fully grown)"
$ws.Range("G4").Value = "5y 28d"
$ws.Range("H4").Value = "Available"
$ws.Range("I4").Value = 18.9

# --- Row 5: Polly -> FALCOR (Dog Holding E, Hold - Surgery) ---
$ws.Range("A5").Value = "FALCOR"
$ws.Range("B5").Value = "58997270"
$ws.Range("C5").Value = "Dog"
$ws.Range("D5").Value = "Dog Holding E"
$ws.Range("E5").Value = "07/30/2025"
$ws.Range("F5").Value = "Mixed Breed, Large (over 44 lbs fully grown)"
$ws.Range("G5").Value = "2y 1d"
$ws.Range("H5").Value = "Hold - Surgery"
$ws.Range("I5").Value = 19.9

# --- Row 6: PAISLEE -> Graffiti (Foster Home, In Foster) ---
$ws.Range("A6").Value = "Graffiti"
$ws.Range("B6").Value = "58834498"
$ws.Range("C6").Value = "Dog"
$ws.Range("D6").Value = "Foster Home"
$ws.Range("E6").Value = "07/02/2025"
$ws.Range("F6").Value = "Bulldog"
$ws.Range("G6").Value = "1y 1m 16d"
$ws.Range("H6").Value = "In Foster"
$ws.Range("I6").Value = 47.6

# --- Row 7: Evelyn -> (blank name) Farm Type Fowl / Chicken, Evaluate ---
$ws.Range("A7").ClearContents()
$ws.Range("B7").Value = "59024880"
$ws.Range("C7").Value = "Farm Type Fowl"
$ws.Range("D7").Value = "Farm"
$ws.Range("E7").Value = "08/01/2025"
$ws.Range("F7").Value = "Chicken"
$ws.Range("G7").ClearContents()
$ws.Range("H7").Value = "Evaluate"
$ws.Range("I7").Value = 17.9

# --- Row 8 (new row): another Farm Type Fowl / Chicken, Evaluate ---
$ws.Range("A8").ClearContents()
$ws.Range("B8").Value = "59024884"
$ws.Range("C8").Value = "Farm Type Fowl"
$ws.Range("D8").Value = "Farm"
$ws.Range("E8").Value = "08/01/2025"
$ws.Range("F8").Value = "Chicken"
$ws.Range("G8").ClearContents()
$ws.Range("H8").Value = "Evaluate"
$ws.Range("I8").Value = 17.9

# --- Row 9 (new row): KAHLUA, Miniature Pig, Farm, Evaluate ---
$ws.Range("A9").Value = "KAHLUA"
$ws.Range("B9").Value = "58972687"
$ws.Range("C9").Value = "Miniature Pig"
$ws.Range("D9").Value = "Farm"
$ws.Range("E9").Value = "07/29/2025"
$ws.Range("F9").Value = "Kune Kune Mini-pig"
$ws.Range("G9").Value = "1y 26d"
$ws.Range("H9").Value = "Evaluate"
$ws.Range("I9").Value = 20.9

# --- Remove the temporary text format now that the values are locked in as text ---
$ws.Range("B2:B9").ClearFormats()
$ws.Range("E2:E9").ClearFormats()

# --- Re-apply AutoFilter over the new extended range A1:I9 ---
$ws.AutoFilterMode = $false
$ws.Range("A1:I9").AutoFilter()

# --- Keep the _FilterDatabase defined name in sync with the new range ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$I`$9"
    }
}

# --- Update column widths to best-fit the new content ---
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(8).AutoFit()

# --- Update print orientation to portrait ---
$ws.PageSetup.Orientation = 1

# --- Move the active selection to F9, matching the last-edited cell ---
$ws.Range("F9").Select()

Write-Host "Edit applied"
